$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.496.68'
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").Value = '2.373.49'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'310.00"
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").Value = "'104.08"
$ws.Range("E6").Value = '  +2.79%  '
$ws.Range("E7").Value = '  -4.64%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -1.01%  '
$ws.Range("D10").Value = "'35.75"
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("D11").Value = "'53.24"
$ws.Range("E11").Value = '  +1.90%  '
$ws.Range("D12").Value = "'0.0806"
$ws.Range("E12").Value = '  -1.50%  '
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("E14").Value = '  -3.87%  '
$ws.Range("D15").Value = "'15.56"
$ws.Range("E15").Value = '  +3.46%  '
$ws.Range("D16").Value = '2.373.77'
$ws.Range("E16").Value = '  +2.36%  '
$ws.Range("D17").Value = "'0.811"
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").Value = '43.461.77'
$ws.Range("E18").Value = '  +0.88%  '
$ws.Range("D19").Value = "'6.31"
$ws.Range("E19").Value = '  +2.83%  '
$ws.Range("D20").Value = "'11.91"
$ws.Range("E20").Value = '  -5.39%  '
$ws.Range("D21").Value = '0.0₃0913'
$ws.Range("E21").Value = '  -0.97%  '
$ws.Range("D22").Value = "'68.22"
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("D23").Value = "'240.07"
$ws.Range("E23").Value = '  -0.69%  '
$ws.Range("E24").Value = '  +0.49%  '
$ws.Range("E25").Value = '  -1.54%  '
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("E27").Value = '  +3.41%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = "'2.33"
$ws.Range("E28").Value = '  +9.97%  '
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").Value = "'36.50"
$ws.Range("E29").Value = '  -2.39%  '
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").Value = "'9.47"
$ws.Range("E30").Value = '  -2.24%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = "'161.05"
$ws.Range("E31").Value = '  -3.08%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = "'5.21"
$ws.Range("E32").Value = '  -2.62%  '
$ws.Range("B33").Value = 'FirstDigitalUSD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("B34").Value = 'Celestia'
$ws.Range("C34").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D34").Value = "'18.19"
$ws.Range("E34").Value = '  +1.11%  '
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = "'2.52"
$ws.Range("E35").Value = '  +5.16%  '
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").Value = "'3.07"
$ws.Range("E36").Value = '  -2.33%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = "'4.65"
$ws.Range("E37").Value = '  +7.66%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = "'0.0734"
$ws.Range("E38").Value = '  -1.21%  '
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").Value = "'1.92"
$ws.Range("E39").Value = '  +3.61%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = "'0.105"
$ws.Range("E40").Value = '  -2.40%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = "'0.114"
$ws.Range("E41").Value = '  -2.19%  '
$ws.Range("B42").Value = 'ApeXProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D42").Value = "'2.61"
$ws.Range("E42").Value = '  +12.92%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.032.02'
$ws.Range("E43").Value = '  +2.60%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = "'19.63"
$ws.Range("E44").Value = '  +0.45%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = "'0.0289"
$ws.Range("E45").Value = '  -0.71%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = "'10.54"
$ws.Range("E46").Value = '  +6.84%  '
$ws.Range("D47").Value = "'3.10"
$ws.Range("E47").Value = '  +2.40%  '
$ws.Range("B48").Value = 'MultiversX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D48").Value = "'57.68"
$ws.Range("E48").Value = '  +3.49%  '
$ws.Range("B49").Value = 'HuobiToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D49").Value = "'2.95"
$ws.Range("E49").Value = '  -0.79%  '
$ws.Range("D50").Value = '2.607.70'
$ws.Range("E50").Value = '  +2.52%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").Value = "'1.56"
$ws.Range("E51").Value = '  +1.33%  '
